$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Update the Create Account test data values
$ws.Range("C2").Value = "abc1234"
$ws.Range("D2").Value = "abc1234"
$ws.Range("B2").Value = "abc2@gmail.com"

# Move the active selection to B2
$ws.Activate()
$ws.Range("B2").Select()
